$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 into row 5 (style, borders, etc.)
$ws.Range("A4:AT4").Copy($ws.Range("A5:AT5"))

# Set the year label for the new row
$ws.Range("A5").Value = "2021年"

# Populate the data values for 2021
$ws.Range("B5").Value = 40.436
$ws.Range("C5").Value = 52.732
$ws.Range("D5").Value = 60.952
$ws.Range("E5").Value = 45.046
$ws.Range("F5").Value = 40
$ws.Range("G5").Value = 39.945
$ws.Range("H5").Value = 41.764
$ws.Range("I5").Value = 44.972
$ws.Range("J5").Value = 38.484
$ws.Range("K5").Value = 58.417
$ws.Range("L5").Value = 39.742
$ws.Range("M5").Value = 40.531
$ws.Range("N5").Value = 29.928
$ws.Range("O5").Value = 29.874
$ws.Range("P5").Value = 38.991
$ws.Range("Q5").Value = 36.102
$ws.Range("R5").Value = 23.444
$ws.Range("S5").Value = 24.035
$ws.Range("T5").Value = 40.343
$ws.Range("U5").Value = 21.019
$ws.Range("V5").Value = 44.649
$ws.Range("W5").Value = 61.111
$ws.Range("X5").Value = 19.479
$ws.Range("Y5").Value = 26.97
$ws.Range("Z5").Value = 21.679
$ws.Range("AA5").Value = 20.422
$ws.Range("AB5").Value = 50.736
$ws.Range("AC5").Value = 27.327
$ws.Range("AD5").Value = 36.811
$ws.Range("AE5").Value = 34.266
$ws.Range("AF5").Value = 31.676
$ws.Range("AG5").Value = 26.451
$ws.Range("AH5").Value = 55.771
$ws.Range("AI5").Value = 47.956
$ws.Range("AJ5").Value = 36.341
$ws.Range("AK5").Value = 48.8
$ws.Range("AL5").Value = 20.74
$ws.Range("AM5").Value = 32.759
$ws.Range("AN5").Value = 38.763
$ws.Range("AO5").Value = 45.482
$ws.Range("AP5").Value = 31.502
$ws.Range("AQ5").Value = 21.014
$ws.Range("AR5").Value = 50.84
$ws.Range("AS5").Value = 30.953
$ws.Range("AT5").Value = 18.453
